$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New vocabulary entries added 2020-12-02 (rows 137-162).
# Column order: A = Chinese sentence (FOREIGN), B = English translation (ENGLISH),
# C = Chinese headword (WORD), D = date added (DATE).

# Row 137
$ws.Cells.Item(137, 1).Value = '存在即被感知。  '
$ws.Cells.Item(137, 2).Value = 'Being is being perceived.'
$ws.Cells.Item(137, 3).Value = '即'
$ws.Cells.Item(137, 4).NumberFormat = "@"
$ws.Cells.Item(137, 4).Value = '2020-12-02'
$ws.Cells.Item(137, 4).Style = "Normal"

# Row 138
$ws.Cells.Item(138, 1).Value = '我打个电话即到。  '
$ws.Cells.Item(138, 2).Value = 'I will be there immediately – I just have to make a phone call.'
$ws.Cells.Item(138, 3).Value = '即'
$ws.Cells.Item(138, 4).NumberFormat = "@"
$ws.Cells.Item(138, 4).Value = '2020-12-02'
$ws.Cells.Item(138, 4).Style = "Normal"

# Row 139
$ws.Cells.Item(139, 1).Value = '康熙一继位即命人修史。  '
$ws.Cells.Item(139, 2).Value = 'Kang Xi ordered people to compile historical books the moment he succeeded to the throne.'
$ws.Cells.Item(139, 3).Value = '即'
$ws.Cells.Item(139, 4).NumberFormat = "@"
$ws.Cells.Item(139, 4).Value = '2020-12-02'
$ws.Cells.Item(139, 4).Style = "Normal"

# Row 140
$ws.Cells.Item(140, 1).Value = '这个化学反应产生很多热量。  '
$ws.Cells.Item(140, 2).Value = 'This chemical reaction generates a lot of heat.'
$ws.Cells.Item(140, 3).Value = '热量'
$ws.Cells.Item(140, 4).NumberFormat = "@"
$ws.Cells.Item(140, 4).Value = '2020-12-02'
$ws.Cells.Item(140, 4).Style = "Normal"

# Row 141
$ws.Cells.Item(141, 1).Value = '气体之间通过热对流传导热量。  '
$ws.Cells.Item(141, 2).Value = 'Gases transmit heat through thermal convention.'
$ws.Cells.Item(141, 3).Value = '热量'
$ws.Cells.Item(141, 4).NumberFormat = "@"
$ws.Cells.Item(141, 4).Value = '2020-12-02'
$ws.Cells.Item(141, 4).Style = "Normal"

# Row 142
$ws.Cells.Item(142, 1).Value = '热力学定律解释了热量转化的原理。  '
$ws.Cells.Item(142, 2).Value = 'The Law of Thermodynamics explains the transfer of heat.'
$ws.Cells.Item(142, 3).Value = '热量'
$ws.Cells.Item(142, 4).NumberFormat = "@"
$ws.Cells.Item(142, 4).Value = '2020-12-02'
$ws.Cells.Item(142, 4).Style = "Normal"

# Row 143
$ws.Cells.Item(143, 1).Value = '要想身体健康，就要多注意饮食营养配比。  '
$ws.Cells.Item(143, 2).Value = 'Pay attention to the amount of nutrition in your food if you want to be healthy.'
$ws.Cells.Item(143, 3).Value = '配比'
$ws.Cells.Item(143, 4).NumberFormat = "@"
$ws.Cells.Item(143, 4).Value = '2020-12-02'
$ws.Cells.Item(143, 4).Style = "Normal"

# Row 144
$ws.Cells.Item(144, 1).Value = '生石灰配比  '
$ws.Cells.Item(144, 2).Value = 'lime proportion.'
$ws.Cells.Item(144, 3).Value = '配比'
$ws.Cells.Item(144, 4).NumberFormat = "@"
$ws.Cells.Item(144, 4).Value = '2020-12-02'
$ws.Cells.Item(144, 4).Style = "Normal"

# Row 145
$ws.Cells.Item(145, 1).Value = '贵金属配比  '
$ws.Cells.Item(145, 2).Value = 'Precious metals ratio.'
$ws.Cells.Item(145, 3).Value = '配比'
$ws.Cells.Item(145, 4).NumberFormat = "@"
$ws.Cells.Item(145, 4).Value = '2020-12-02'
$ws.Cells.Item(145, 4).Style = "Normal"

# Row 146
$ws.Cells.Item(146, 1).Value = '桌子摆放得很起眼儿。  '
$ws.Cells.Item(146, 2).Value = 'The table was set appealingly.'
$ws.Cells.Item(146, 3).Value = '摆放'
$ws.Cells.Item(146, 4).NumberFormat = "@"
$ws.Cells.Item(146, 4).Value = '2020-12-02'
$ws.Cells.Item(146, 4).Style = "Normal"

# Row 147
$ws.Cells.Item(147, 1).Value = '茶几摆放得非常合适。  '
$ws.Cells.Item(147, 2).Value = 'End tables placed conveniently.'
$ws.Cells.Item(147, 3).Value = '摆放'
$ws.Cells.Item(147, 4).NumberFormat = "@"
$ws.Cells.Item(147, 4).Value = '2020-12-02'
$ws.Cells.Item(147, 4).Style = "Normal"

# Row 148
$ws.Cells.Item(148, 1).Value = '这些家具都毫无用处地摆放着。  '
$ws.Cells.Item(148, 2).Value = 'The furniture was sitting around uselessly.'
$ws.Cells.Item(148, 3).Value = '摆放'
$ws.Cells.Item(148, 4).NumberFormat = "@"
$ws.Cells.Item(148, 4).Value = '2020-12-02'
$ws.Cells.Item(148, 4).Style = "Normal"

# Row 149
$ws.Cells.Item(149, 1).Value = '他们肯定要摄入大量的水份。  '
$ws.Cells.Item(149, 2).Value = 'They will certainly need to take in plenty of liquid.'
$ws.Cells.Item(149, 3).Value = '摄入'
$ws.Cells.Item(149, 4).NumberFormat = "@"
$ws.Cells.Item(149, 4).Value = '2020-12-02'
$ws.Cells.Item(149, 4).Style = "Normal"

# Row 150
$ws.Cells.Item(150, 1).Value = '在健康状态下，液体的摄入与排出应是均衡的。  '
$ws.Cells.Item(150, 2).Value = 'In health，fluid intake is usually balanced with output.'
$ws.Cells.Item(150, 3).Value = '摄入'
$ws.Cells.Item(150, 4).NumberFormat = "@"
$ws.Cells.Item(150, 4).Value = '2020-12-02'
$ws.Cells.Item(150, 4).Style = "Normal"

# Row 151
$ws.Cells.Item(151, 1).Value = '表4是一个典型的一天食物摄入的范例。  '
$ws.Cells.Item(151, 2).Value = 'Table 4 is an example of a typical day''s food intake.'
$ws.Cells.Item(151, 3).Value = '摄入'
$ws.Cells.Item(151, 4).NumberFormat = "@"
$ws.Cells.Item(151, 4).Value = '2020-12-02'
$ws.Cells.Item(151, 4).Style = "Normal"

# Row 152
$ws.Cells.Item(152, 1).Value = '腹鳍对应着四足动物的后肢。  '
$ws.Cells.Item(152, 2).Value = 'Ventral (or pelvic) fins correspond to the hind limbs of a quadruped.'
$ws.Cells.Item(152, 3).Value = '对应'
$ws.Cells.Item(152, 4).NumberFormat = "@"
$ws.Cells.Item(152, 4).Value = '2020-12-02'
$ws.Cells.Item(152, 4).Style = "Normal"

# Row 153
$ws.Cells.Item(153, 1).Value = '蝙蝠的翅膀和人的胳膊是对应的。  '
$ws.Cells.Item(153, 2).Value = 'The wing of a bat and the arm of a man are homologous.'
$ws.Cells.Item(153, 3).Value = '对应'
$ws.Cells.Item(153, 4).NumberFormat = "@"
$ws.Cells.Item(153, 4).Value = '2020-12-02'
$ws.Cells.Item(153, 4).Style = "Normal"

# Row 154
$ws.Cells.Item(154, 1).Value = '我们可以发现这些抗体是特定对应每种蛋白质的。  '
$ws.Cells.Item(154, 2).Value = 'We were able to see that these antibodies are specific to each protein.'
$ws.Cells.Item(154, 3).Value = '对应'
$ws.Cells.Item(154, 4).NumberFormat = "@"
$ws.Cells.Item(154, 4).Value = '2020-12-02'
$ws.Cells.Item(154, 4).Style = "Normal"

# Row 155
$ws.Cells.Item(155, 1).Value = '除此以外  '
$ws.Cells.Item(155, 2).Value = 'Put that aside.'
$ws.Cells.Item(155, 3).Value = '除此以外'
$ws.Cells.Item(155, 4).NumberFormat = "@"
$ws.Cells.Item(155, 4).Value = '2020-12-02'
$ws.Cells.Item(155, 4).Style = "Normal"

# Row 156
$ws.Cells.Item(156, 1).Value = '好除此以外  '
$ws.Cells.Item(156, 2).Value = 'Yes, put that aside.'
$ws.Cells.Item(156, 3).Value = '除此以外'
$ws.Cells.Item(156, 4).NumberFormat = "@"
$ws.Cells.Item(156, 4).Value = '2020-12-02'
$ws.Cells.Item(156, 4).Style = "Normal"

# Row 157
$ws.Cells.Item(157, 1).Value = '除此以外，看不出还有什么别的原因会使她与达切发生一段短暂的、充满不安的恋情。  '
$ws.Cells.Item(157, 2).Value = 'No other explanation is offered for her decision to enter into a brief, nervous affair with Dutch.'
$ws.Cells.Item(157, 3).Value = '除此以外'
$ws.Cells.Item(157, 4).NumberFormat = "@"
$ws.Cells.Item(157, 4).Value = '2020-12-02'
$ws.Cells.Item(157, 4).Style = "Normal"

# Row 158
$ws.Cells.Item(158, 1).Value = '利益逐渐消失。  '
$ws.Cells.Item(158, 2).Value = 'Interest tapered off.'
$ws.Cells.Item(158, 3).Value = '逐渐'
$ws.Cells.Item(158, 4).NumberFormat = "@"
$ws.Cells.Item(158, 4).Value = '2020-12-02'
$ws.Cells.Item(158, 4).Style = "Normal"

# Row 159
$ws.Cells.Item(159, 1).Value = '风逐渐停了。  '
$ws.Cells.Item(159, 2).Value = 'The wind gradually died away.'
$ws.Cells.Item(159, 3).Value = '逐渐'
$ws.Cells.Item(159, 4).NumberFormat = "@"
$ws.Cells.Item(159, 4).Value = '2020-12-02'
$ws.Cells.Item(159, 4).Style = "Normal"

# Row 160
$ws.Cells.Item(160, 1).Value = '利润逐渐下滑。  '
$ws.Cells.Item(160, 2).Value = 'The interest declined by little and little.'
$ws.Cells.Item(160, 3).Value = '逐渐'
$ws.Cells.Item(160, 4).NumberFormat = "@"
$ws.Cells.Item(160, 4).Value = '2020-12-02'
$ws.Cells.Item(160, 4).Style = "Normal"

# Row 161
$ws.Cells.Item(161, 1).Value = '习惯是后天养成的。  '
$ws.Cells.Item(161, 2).Value = 'Xíguàn shìhòutiān yǎng chéng de. '
$ws.Cells.Item(161, 3).Value = '养成'
$ws.Cells.Item(161, 4).NumberFormat = "@"
$ws.Cells.Item(161, 4).Value = '2020-12-02'
$ws.Cells.Item(161, 4).Style = "Normal"

# Row 162
$ws.Cells.Item(162, 2).Value = 'Yǎngchéng jiàoyù hěn zhòngyào，yào zhùyì péiyǎng háizi de liánghǎo xíguàn.'
$ws.Cells.Item(162, 3).Value = '养成'
$ws.Cells.Item(162, 4).NumberFormat = "@"
$ws.Cells.Item(162, 4).Value = '2020-12-02'
$ws.Cells.Item(162, 4).Style = "Normal"

